$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 122.916664
$ws.Range("I9").Value = 136.8
$ws.Range("J9").Value = 53.5
$ws.Range("K9").Value = 136.8
$ws.Range("L9").Value = 53.5
$ws.Range("M9").Value = 32.19999999999999
$ws.Range("N9").Value = -391.5
$ws.Range("H15").Value = 1208.0385
$ws.Range("I15").Value = 1208.0385
$ws.Range("K15").Value = 3624.1155
$ws.Range("M15").Value = -3455.1155
$ws.Range("H17").Value = 635.6
$ws.Range("J17").Value = 688.94116
$ws.Range("L17").Value = 2066.82348
$ws.Range("N17").Value = -2402.82348
$ws.Range("H40").Value = 13883.857
$ws.Range("I40").Value = 3882.8
$ws.Range("J40").Value = 17009.188
$ws.Range("K40").Value = 3882.8
$ws.Range("L40").Value = 17009.188
$ws.Range("M40").Value = -3707.8
$ws.Range("N40").Value = -17359.188
$ws.Range("H101").Value = 713
$ws.Range("I101").Value = 728.3333
$ws.Range("J101").Value = 690
$ws.Range("K101").Value = 2184.9999
$ws.Range("L101").Value = 2070
$ws.Range("M101").Value = -562.9998999999998
$ws.Range("N101").Value = -5314
$ws.Range("H121").Value = 2395.6875
$ws.Range("J121").Value = 2395.6875
$ws.Range("L121").Value = 7187.0625
$ws.Range("N121").Value = -10681.0625
$ws.Range("H127").Value = 1882
$ws.Range("I127").Value = 1475.4
$ws.Range("J127").Value = 2898.5
$ws.Range("K127").Value = 4426.200000000001
$ws.Range("L127").Value = 8695.5
$ws.Range("M127").Value = 533.7999999999993
$ws.Range("N127").Value = -18615.5
$ws.Range("H135").Value = 21621
$ws.Range("I135").Value = 2000
$ws.Range("K135").Value = 18000
$ws.Range("M135").Value = -15465
$ws.Range("H138").Value = 5504.283
$ws.Range("J138").Value = 6857.657
$ws.Range("L138").Value = 20572.971
$ws.Range("N138").Value = -30852.971
$ws.Range("H140").Value = 86780
$ws.Range("J140").Value = 86780
$ws.Range("L140").Value = 86780
$ws.Range("N140").Value = -97140
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2456.946
$ws.Range("I2").Value = 2085.6206
$ws.Range("K2").Value = 2085.6206
$ws.Range("M2").Value = -1972.6206
$ws.Range("H32").Value = 5222.162
$ws.Range("I32").Value = 456.51562
$ws.Range("K32").Value = 456.51562
$ws.Range("M32").Value = -169.51562
$ws.Range("H61").Value = 6219.636
$ws.Range("I61").Value = 6491.6
$ws.Range("K61").Value = 6491.6
$ws.Range("M61").Value = -6279.6
$ws.Range("H110").Value = 1689.6522
$ws.Range("I110").Value = 1654.15
$ws.Range("J110").Value = 1926.3334
$ws.Range("K110").Value = 1654.15
$ws.Range("L110").Value = 1926.3334
$ws.Range("M110").Value = 390.8499999999999
$ws.Range("N110").Value = -6016.3334
$ws.Range("H116").Value = 2456.946
$ws.Range("I116").Value = 2085.6206
$ws.Range("K116").Value = 2085.6206
$ws.Range("M116").Value = 208.3793999999998
$ws.Range("H132").Value = 4449.4707
$ws.Range("I132").Value = 4415.125
$ws.Range("K132").Value = 13245.375
$ws.Range("M132").Value = -10715.375
$ws.Range("H136").Value = 6219.636
$ws.Range("I136").Value = 6491.6
$ws.Range("K136").Value = 19474.8
$ws.Range("M136").Value = -16924.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2456.946
$ws.Range("I3").Value = 2085.6206
$ws.Range("K3").Value = 2085.6206
$ws.Range("M3").Value = -1971.6206
$ws.Range("H86").Value = 2755.606
$ws.Range("I86").Value = 2998.087
$ws.Range("J86").Value = 2197.9
$ws.Range("K86").Value = 2998.087
$ws.Range("L86").Value = 2197.9
$ws.Range("M86").Value = -1875.087
$ws.Range("N86").Value = -4443.9
$ws.Range("H89").Value = 2755.606
$ws.Range("I89").Value = 2998.087
$ws.Range("J89").Value = 2197.9
$ws.Range("K89").Value = 14990.435
$ws.Range("L89").Value = 10989.5
$ws.Range("M89").Value = -9374.434999999999
$ws.Range("N89").Value = -22221.5
$ws.Range("H125").Value = 139997
$ws.Range("J125").Value = 139997
$ws.Range("L125").Value = 139997
$ws.Range("N125").Value = -149837
$ws.Range("H126").Value = 49999.332
$ws.Range("J126").Value = 57499
$ws.Range("L126").Value = 57499
$ws.Range("N126").Value = -67379
$ws.Range("H134").Value = 2078.4
$ws.Range("I134").Value = 1598
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 4794
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -2259
$ws.Range("N134").Value = -17070
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5034.3335
$ws.Range("I58").Value = 5307.778
$ws.Range("J58").Value = 4214
$ws.Range("K58").Value = 5307.778
$ws.Range("L58").Value = 4214
$ws.Range("M58").Value = -5104.778
$ws.Range("N58").Value = -4620
$ws.Range("H86").Value = 9114
$ws.Range("J86").Value = 9449.75
$ws.Range("L86").Value = 9449.75
$ws.Range("N86").Value = -11695.75
$ws.Range("H89").Value = 9114
$ws.Range("J89").Value = 9449.75
$ws.Range("L89").Value = 47248.75
$ws.Range("N89").Value = -58480.75
$ws.Range("H122").Value = 3739.353
$ws.Range("I122").Value = 3938.182
$ws.Range("J122").Value = 3374.8333
$ws.Range("K122").Value = 11814.546
$ws.Range("L122").Value = 10124.4999
$ws.Range("M122").Value = -9364.545999999998
$ws.Range("N122").Value = -15024.4999
$ws.Range("H134").Value = 3562.9614
$ws.Range("I134").Value = 2586.95
$ws.Range("J134").Value = 6816.3335
$ws.Range("K134").Value = 7760.849999999999
$ws.Range("L134").Value = 20449.0005
$ws.Range("M134").Value = -5225.849999999999
$ws.Range("N134").Value = -25519.0005
$ws.Range("H136").Value = 5034.3335
$ws.Range("I136").Value = 5307.778
$ws.Range("J136").Value = 4214
$ws.Range("K136").Value = 15923.334
$ws.Range("L136").Value = 12642
$ws.Range("M136").Value = -13373.334
$ws.Range("N136").Value = -17742
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 101038.8
$ws.Range("I46").Value = 438.8
$ws.Range("K46").Value = 1316.4
$ws.Range("M46").Value = -1225.4
$ws.Range("H68").Value = 2170.5264
$ws.Range("J68").Value = 2292.3572
$ws.Range("L68").Value = 6877.071599999999
$ws.Range("N68").Value = -8499.071599999999
$ws.Range("H71").Value = 2170.5264
$ws.Range("J71").Value = 2292.3572
$ws.Range("L71").Value = 20631.2148
$ws.Range("N71").Value = -28743.2148
$ws.Range("H116").Value = 8876.333000000001
$ws.Range("I116").Value = 8876.333000000001
$ws.Range("K116").Value = 26628.999
$ws.Range("M116").Value = -23186.999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 22740.666
$ws.Range("I55").Value = 22740.666
$ws.Range("K55").Value = 22740.666
$ws.Range("M55").Value = -22413.666
$ws.Range("H132").Value = 4763
$ws.Range("I132").Value = 6158.25
$ws.Range("J132").Value = 3832.8333
$ws.Range("K132").Value = 18474.75
$ws.Range("L132").Value = 11498.4999
$ws.Range("M132").Value = -15944.75
$ws.Range("N132").Value = -16558.4999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4922.6304
$ws.Range("I40").Value = 3069.516
$ws.Range("J40").Value = 8752.4
$ws.Range("K40").Value = 3069.516
$ws.Range("L40").Value = 8752.4
$ws.Range("M40").Value = -2933.516
$ws.Range("N40").Value = -9024.4
$ws.Range("H46").Value = 2179
$ws.Range("I46").Value = 1250.5
$ws.Range("K46").Value = 1250.5
$ws.Range("M46").Value = -1062.5
$ws.Range("H125").Value = 91854.57000000001
$ws.Range("J125").Value = 91854.57000000001
$ws.Range("L125").Value = 91854.57000000001
$ws.Range("N125").Value = -101694.57
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1089.3636
$ws.Range("I96").Value = 1017.1667
$ws.Range("J96").Value = 1176
$ws.Range("K96").Value = 1017.1667
$ws.Range("L96").Value = 1176
$ws.Range("M96").Value = 355.8333
$ws.Range("N96").Value = -3922
$ws.Range("H126").Value = 1952.1333
$ws.Range("I126").Value = 1744.2727
$ws.Range("J126").Value = 2523.75
$ws.Range("K126").Value = 5232.8181
$ws.Range("L126").Value = 7571.25
$ws.Range("M126").Value = -2762.8181
$ws.Range("N126").Value = -12511.25
